# Generate Report for Handoff
# Rename the source file GUID from 42ba1bf8-8419-470d-bd34-dd6d4dc298e3 to
# 94e35117-715b-43b5-9d4e-54bc1dad67bb across all three sheets, refresh the
# handoff xliff file names (new content hash), and bump the handoff/generate
# timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "42ba1bf8-8419-470d-bd34-dd6d4dc298e3"
$newGuid = "94e35117-715b-43b5-9d4e-54bc1dad67bb"

$oldHash = "c8a3820de00f016c990393f5c2d65c348894ee7e"
$newHash = "016a76dd21889117c16de60df1eb254461145ebb"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b163ef8833b9cbc6457089651bb94ac8b4bd0ef/e2e/"

# The sandboxed host re-materialises a cell's style when a hyperlink is
# (re)built via Hyperlinks.Add, swapping the workbook's custom "HyperLink"
# font for a generic one. Restore the original look (single underline,
# cornflower-blue text) right after, so the visible formatting matches.
$hyperlinkColor = 15570276  # OLE BGR for RGB(0x64,0x95,0xED) == style "FF6495ED"

function Set-CellHyperlink($ws, $cellAddr, $address, $displayText) {
    $ws.Range($cellAddr).Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $address, "", "", $displayText)
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid.md"

Set-CellHyperlink $ws "B2" "$repoBase$newGuid.md" "e2e\$newGuid.md"

$ws.Range("G2").Value = "2016-08-30 23:02:56"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

Set-CellHyperlink $ws "A2" "$repoBase$newGuid.md" "$newGuid.md"

$ws.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-30 23:02:52"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

Set-CellHyperlink $ws "A2" "$repoBase$newGuid.md" "$newGuid.md"

$ws.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-30 23:02:56"
